# Update "想去人数" (interest count) figures that changed between crawls.
# Same underlying events are listed both on their category sheet
# (展览 / 演出 / 本地生活) and on the combined "全部类型" sheet, so each
# updated figure is written to both places.

$wb = $excel.ActiveWorkbook

# 展览
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    "F2"  = 13794
    "F6"  = 808
    "F8"  = 212
    "F9"  = 140
    "F10" = 123
    "F11" = 257
    "F13" = 625
    "F15" = 547
    "F17" = 39
    "F18" = 323
    "F20" = 171
    "F21" = 97
    "F22" = 52
    "F25" = 127
    "F26" = 49
}
foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# 演出
$wsShow = $wb.Worksheets.Item("演出")
$showUpdates = @{
    "F6"  = 143
    "F7"  = 187
    "F8"  = 2174
    "F15" = 1932
}
foreach ($cell in $showUpdates.Keys) {
    $wsShow.Range($cell).Value = $showUpdates[$cell]
}

# 本地生活
$wsLocal = $wb.Worksheets.Item("本地生活")
$localUpdates = @{
    "F2" = 241
    "F3" = 220
}
foreach ($cell in $localUpdates.Keys) {
    $wsLocal.Range($cell).Value = $localUpdates[$cell]
}

# 全部类型 (combined sheet containing the same records)
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    "F2"  = 241
    "F3"  = 13794
    "F7"  = 808
    "F11" = 220
    "F12" = 212
    "F13" = 140
    "F14" = 123
    "F15" = 257
    "F19" = 143
    "F21" = 625
    "F23" = 547
    "F25" = 39
    "F26" = 323
    "F28" = 187
    "F29" = 2174
    "F34" = 171
    "F35" = 97
    "F36" = 52
    "F41" = 127
    "F42" = 49
    "F43" = 1932
}
foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
